# Insert 3 new "guarda" rows for Camote (Zapallo) before the existing
# row 48, pushing the old rows 48-83 down to 51-86, then fill in the
# new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 48..83 down by 3 (creates 3 blank rows at 48:50, carrying
# the existing row formatting, e.g. the date number format in column D).
$ws.Rows("48:50").Insert()

# Common / constant column values for this Camote (Zapallo) block.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$catId     = 100112045
$categoria = "Zapallo"
$variedad  = "Camote"
$unidad    = "`$/kilo (volumen en unidades)"
$kgUnid    = 1
$clasif    = "Hortaliza"
$fecha     = 45089

# Row 48: 1a (guarda)
$ws.Range("A48").Value = $mercadoId
$ws.Range("B48").Value = $mercado
$ws.Range("C48").Value = $region
$ws.Range("D48").Value = $fecha
$ws.Range("E48").Value = $codreg
$ws.Range("F48").Value = $catId
$ws.Range("G48").Value = $categoria
$ws.Range("H48").Value = $variedad
$ws.Range("I48").Value = "1a (guarda)"
$ws.Range("J48").Value = 400
$ws.Range("K48").Value = 400
$ws.Range("L48").Value = 450
$ws.Range("M48").Value = 425
$ws.Range("N48").Value = $unidad
$ws.Range("O48").Value = "Región de O'Higgins"
$ws.Range("P48").Value = 425
$ws.Range("Q48").Value = $kgUnid
$ws.Range("R48").Value = $clasif

# Row 49: 2a (guarda)
$ws.Range("A49").Value = $mercadoId
$ws.Range("B49").Value = $mercado
$ws.Range("C49").Value = $region
$ws.Range("D49").Value = $fecha
$ws.Range("E49").Value = $codreg
$ws.Range("F49").Value = $catId
$ws.Range("G49").Value = $categoria
$ws.Range("H49").Value = $variedad
$ws.Range("I49").Value = "2a (guarda)"
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 380
$ws.Range("L49").Value = 400
$ws.Range("M49").Value = 390
$ws.Range("N49").Value = $unidad
$ws.Range("O49").Value = "Región de O'Higgins"
$ws.Range("P49").Value = 390
$ws.Range("Q49").Value = $kgUnid
$ws.Range("R49").Value = $clasif

# Row 50: 3a (guarda)
$ws.Range("A50").Value = $mercadoId
$ws.Range("B50").Value = $mercado
$ws.Range("C50").Value = $region
$ws.Range("D50").Value = $fecha
$ws.Range("E50").Value = $codreg
$ws.Range("F50").Value = $catId
$ws.Range("G50").Value = $categoria
$ws.Range("H50").Value = $variedad
$ws.Range("I50").Value = "3a (guarda)"
$ws.Range("J50").Value = 400
$ws.Range("K50").Value = 360
$ws.Range("L50").Value = 380
$ws.Range("M50").Value = 370
$ws.Range("N50").Value = $unidad
$ws.Range("O50").Value = "Región de O'Higgins"
$ws.Range("P50").Value = 370
$ws.Range("Q50").Value = $kgUnid
$ws.Range("R50").Value = $clasif
